$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 0.04554333333333333
$ws.Range("N2").Value = 0.13663
$ws.Range("Q2").Value = 0.04810358217888889
$ws.Range("R2").Value = 0.43293223961
